# Replace commas with full stops in the size (K) column.
# Cells that currently hold text values "1,5", "2,5" or "3,5" should
# become real numbers 1.5, 2.5 and 3.5 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOT")

$map = @{ "1,5" = 1.5; "2,5" = 2.5; "3,5" = 3.5 }

$used = $ws.Range("K1:K999")
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($val -ne $null -and $map.ContainsKey([string]$val)) {
        $cell.Value = $map[[string]$val]
    }
}

# Update the sheet view/selection like the recorded session: select
# the whole K column and reset the scroll position.
$ws.Activate()
$ws.Range("K1:K1048576").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
